# 25 Mayis verileri eklendi (May 25 data added)
# Adds a new row (75) to the "data" worksheet / Table3 with that day's
# test / case / death / recovered figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row: date (2020-05-25 as serial 43976), test, case, death, recovered
$ws.Range("A75").Value = 43976
$ws.Range("B75").Value = 21492
$ws.Range("C75").Value = 987
$ws.Range("D75").Value = 29
$ws.Range("E75").Value = 1321

# Grow the worksheet table (Table3) so its range / autofilter cover the new row
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E75"))

# Match the author's recorded selection after adding the row
$ws.Range("E74").Select()
